$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (FAPs -> MuSCs target-cluster pairing) was removed entirely from the
# source data re-run; delete it so the sheet only keeps the FAPs -> ECs row.
$ws.Rows.Item(3).Delete()

# The remaining row's derived-specificity / edge-weight columns were
# recomputed against the new TPM values.
$ws.Range("M2").Value = 0.1055746666666667
$ws.Range("N2").Value = 0.316724
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.2211448612408889
$ws.Range("R2").Value = 1.990303751168
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1
